$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Group and Self Assessment")

$ws.Range("D13").Value = 5
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = 5
$ws.Range("G13").Value = 4
